$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered cell => new value updates (from diff).
$updates = [ordered]@{
    'D2' = '92.169.55'
    'E2' = '  +1.76%  '
    'D3' = '3.107.57'
    'E3' = '  -2.42%  '
    'E4' = '  -0.13%  '
    'D5' = '237.18'
    'E5' = '  -1.06%  '
    'D6' = '615.98'
    'E6' = '  -0.37%  '
    'E7' = '  -2.15%  '
    'D8' = '0.390'
    'E8' = '  +4.64%  '
    'E9' = '  -0.03%  '
    'D10' = '3.105.49'
    'E10' = '  -2.15%  '
    'E11' = '  -0.46%  '
    'D12' = '0.200'
    'E12' = '  -1.14%  '
    'D13' = '0.0000248'
    'E13' = '  +0.25%  '
    'D14' = '92.264.23'
    'D15' = '34.25'
    'E15' = '  -2.61%  '
    'D16' = '5.44'
    'E16' = '  -2.20%  '
    'E17' = '  -1.26%  '
    'E18' = '  -1.62%  '
    'E19' = '  +0.22%  '
    'D20' = '14.64'
    'E20' = '  -3.33%  '
    'D21' = '5.77'
    'E21' = '  -4.11%  '
    'D22' = '9.38'
    'E22' = '  +2.31%  '
    'D23' = '446.10'
    'E23' = '  -1.08%  '
    'D24' = '0.0000196'
    'E24' = '  -3.19%  '
    'D25' = '5.77'
    'E25' = '  +0.47%  '
    'D26' = '86.66'
    'E26' = '  +4.90%  '
    'D27' = '11.84'
    'E27' = '  -1.14%  '
    'D28' = '3.275.38'
    'E28' = '  -1.67%  '
    'E29' = '  -0.05%  '
    'D30' = '0.133'
    'E30' = '  -5.86%  '
    'D31' = '0.233'
    'E31' = '  -1.83%  '
    'E32' = '  -0.86%  '
    'D33' = '9.13'
    'E33' = '  -2.24%  '
    'D34' = '0.993'
    'E34' = '  -0.77%  '
    'D35' = '7.87'
    'E35' = '  +2.48%  '
    'E36' = '  -7.67%  '
    'D37' = '26.14'
    'E37' = '  -2.00%  '
    'E38' = '  -3.73%  '
    'E39' = '  +0.90%  '
    'D40' = '484.69'
    'E40' = '  -4.92%  '
    'E41' = '  -4.05%  '
    'E42' = '  +8.21%  '
    'E43' = '  -3.91%  '
    'D44' = '3.30'
    'E44' = '  -4.28%  '
    'D46' = '162.64'
    'E46' = '  +3.83%  '
    'B47' = 'Stacks'
    'C47' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'D47' = '1.89'
    'E47' = '  -2.23%  '
    'B48' = 'ARBITRUM'
    'C48' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D48' = '0.690'
    'E48' = '  -5.25%  '
    'E49' = '  +0.46%  '
    'E50' = '  +3.82%  '
    'E51' = '  -1.29%  '
}

# Force each target cell to Text format before writing so that
# numeric-looking strings (e.g. "237.18", "0.0000248") are kept as
# literal text instead of being parsed into IEEE754 doubles -
# matching the inlineStr/string cells in the source workbook.
# ClearFormats() afterwards drops the temporary style again so the
# cell keeps its original (default) style index.
foreach ($cell in $updates.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$cell]
    $rng.ClearFormats()
}
